$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 data: name, surname, email, interest
$ws.Range("D2").Value = "stocks"
$ws.Range("A2").Value = "Elliot"
$ws.Range("B2").Value = "Sackman"

# Clear rows 3 and 4 (remove the other two people entirely)
$ws.Range("A3:D4").ClearContents()
$ws.Rows.Item(3).RowHeight = 20
$ws.Rows.Item(4).RowHeight = 14.75

# Remove existing hyperlinks, then re-add the hyperlink on C2 with new address
$ws.Hyperlinks.Delete()
$ws.Range("C2").Value = "pythonprojectemail23@gmail.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:pythonprojectemail23@gmail.com")

# Final selection lands on C3
$ws.Range("C3").Select() | Out-Null
